# Split disease analyses into 2014-onwards and 2015-onwards
#
# Style-sheet changes:
#   1. Add a new "Abstract Title" paragraph style (based on Normal, followed
#      by Abstract) used to head the abstract section.
#   2. Abstract style: tighten the space-before from 300 to 100 (since the
#      new Abstract Title style now supplies the space-before=300 above it).
#   3. Add a new "Footnote Block Text" paragraph style (based on Footnote
#      Text) for block-quoted footnote content.

$d = $word.ActiveDocument

# --- 1. New "Abstract Title" style ---------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles.Item("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles.Item("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# --- 2. "Abstract" style: before-spacing 300 -> 100 -----------------------
$abstract = $d.Styles.Item("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. New "Footnote Block Text" style -----------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = $d.Styles.Item("Footnote Text")
$footnoteBlockText.NextParagraphStyle = $d.Styles.Item("Footnote Text")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24
